# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''41.948.42'
$ws.Range("E2").Value = '  +5.55%  '
$ws.Range("D3").Value = '''2.232.00'
$ws.Range("E3").Value = '  +2.50%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '''231.61'
$ws.Range("E5").Value = '  +2.05%  '
$ws.Range("D6").Value = '''0.629'
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("D7").Value = '''61.68'
$ws.Range("E7").Value = '  -2.24%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +2.90%  '
$ws.Range("D10").Value = '''59.12'
$ws.Range("E10").Value = '  +1.07%  '
$ws.Range("D11").Value = '''0.0895'
$ws.Range("E11").Value = '  +4.88%  '
$ws.Range("E12").Value = '  -0.23%  '
$ws.Range("D13").Value = '''2.561.30'
$ws.Range("E13").Value = '  +2.52%  '
$ws.Range("D14").Value = '''15.67'
$ws.Range("E14").Value = '  -1.67%  '
$ws.Range("D15").Value = '''22.04'
$ws.Range("E15").Value = '  +0.68%  '
$ws.Range("D16").Value = '''0.803'
$ws.Range("E16").Value = '  -1.27%  '
$ws.Range("D18").Value = '''2.255.52'
$ws.Range("E18").Value = '  +3.79%  '
$ws.Range("D19").Value = '''41.838.77'
$ws.Range("E19").Value = '  +5.32%  '
$ws.Range("D20").Value = '''72.15'
$ws.Range("E20").Value = '  +0.45%  '
$ws.Range("E21").Value = '  -2.40%  '
$ws.Range("D22").Value = '''6.04'
$ws.Range("E22").Value = '  +0.71%  '
$ws.Range("D23").Value = '''249.51'
$ws.Range("E23").Value = '  +8.50%  '
$ws.Range("D24").Value = '''1.00'
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("E25").Value = '  +2.23%  '
$ws.Range("E26").Value = '  -0.23%  '
$ws.Range("D27").Value = '''9.68'
$ws.Range("E27").Value = '  +0.63%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '''0.142'
$ws.Range("E28").Value = '  +1.26%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = '''167.05'
$ws.Range("E29").Value = '  -2.33%  '
$ws.Range("D30").Value = '''20.03'
$ws.Range("E30").Value = '  +0.77%  '
$ws.Range("E31").Value = '  -3.05%  '
$ws.Range("E32").Value = '  +1.31%  '
$ws.Range("D33").Value = '''0.123'
$ws.Range("E33").Value = '  +0.23%  '
$ws.Range("E34").Value = '  +7.46%  '
$ws.Range("D35").Value = '''4.67'
$ws.Range("E35").Value = '  +2.88%  '
$ws.Range("E36").Value = '  +3.03%  '
$ws.Range("D37").Value = '''6.64'
$ws.Range("E37").Value = '  -4.89%  '
$ws.Range("D38").Value = '''3.69'
$ws.Range("E38").Value = '  -5.46%  '
$ws.Range("D39").Value = '''2.37'
$ws.Range("E39").Value = '  -1.63%  '
$ws.Range("D40").Value = '''0.000268'
$ws.Range("E40").Value = '  +39.22%  '
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("E42").Value = '  +4.86%  '
$ws.Range("D43").Value = '''4.84'
$ws.Range("E43").Value = '  -3.75%  '
$ws.Range("D44").Value = '''8.56'
$ws.Range("E44").Value = '  +8.33%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").Value = '''0.0975'
$ws.Range("E45").Value = '  +6.12%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").Value = '''1.22'
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").Value = '''99.07'
$ws.Range("E47").Value = '  -3.43%  '
$ws.Range("D48").Value = '''1.480.40'
$ws.Range("E48").Value = '  -2.11%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '''16.51'
$ws.Range("E49").Value = '  -6.91%  '
$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D50").Value = '''2.81'
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").Value = '''52.57'
$ws.Range("E51").Value = '  +5.93%  '
